$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '95.932.45'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.547.72'
$ws.Range("E3").Value = '  -1.26%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '239.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.25%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '649.72'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.36%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.62'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.03%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.403'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.42%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.06'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.12%  '
$ws.Range("E10").Value = '  +0.01%  '
$ws.Range("D11").Value = '3.545.04'
$ws.Range("E11").Value = '  -1.28%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '43.17'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +1.06%  '
$ws.Range("E13").Value = '  +0.66%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.35'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.44%  '
$ws.Range("D15").Value = '4.207.92'
$ws.Range("E15").Value = '  -1.37%  '
$ws.Range("D16").Value = '95.847.30'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000258'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.53%  '
$ws.Range("D18").Value = '3.544.85'
$ws.Range("E18").Value = '  -1.40%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.93'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.39'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.520'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '503.70'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -7.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.86'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000196'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.13%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '95.58'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.61'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.88%  '
$ws.Range("D29").Value = '3.739.68'
$ws.Range("E29").Value = '  -1.58%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.151'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +7.90%  '
$ws.Range("E31").Value = '  -4.24%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.30'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.04%  '
$ws.Range("E34").Value = '  +2.30%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.997'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '31.15'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '608.92'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.42%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.67'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.67%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.560'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("E40").Value = '  +7.74%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.149'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.46%  '
$ws.Range("E43").Value = '  -2.60%  '
$ws.Range("E44").Value = '  +4.59%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.67'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.38%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.50'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.06%  '
$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.26'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.47%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0418'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '33.60'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.05%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '8.11'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.54%  '
